# Pseudonymization / field-addition edit:
# Add a "Spouse" column (E) to Sheet1, populating the two rows whose
# subjects have a well-known on-screen spouse (Hank Hill -> Peggy Hill,
# Homer Simpson -> Marge Simpson), leaving the other two rows blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header for the new column
$ws.Range("E1").Value = "Spouse"

# Known spouses for two of the four people already in the sheet
$ws.Range("E2").Value = "Peggy Hill"
$ws.Range("E5").Value = "Marge Simpson"

# Size the new column similarly to its neighbours (best-fit-style width)
$ws.Columns.Item(5).ColumnWidth = 13.8

# Move the active selection the way the source workbook leaves it
$ws.Range("E4").Select() | Out-Null
